$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.212.20'
$ws.Range("E2").Value = '  +2.40%  '

$ws.Range("D3").Value = '2.525.11'
$ws.Range("E3").Value = '  +1.66%  '

$ws.Range("E4").Value = '  +0.06%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '322.82'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.03%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '109.61'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.59%  '

$ws.Range("E7").Value = '  +2.23%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.02%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.554'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +4.04%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '40.60'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +4.85%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '20.54'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +13.44%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0824'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +2.15%  '

$ws.Range("E13").Value = '  +1.26%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '7.29'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.75%  '

$ws.Range("D15").Value = '2.916.19'
$ws.Range("E15").Value = '  +1.77%  '

$ws.Range("D16").Value = '2.523.81'
$ws.Range("E16").Value = '  +1.82%  '

$ws.Range("E17").Value = '  +1.24%  '

$ws.Range("D18").Value = '48.101.61'
$ws.Range("E18").Value = '  +2.35%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '13.43'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +5.90%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.65'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.87%  '

$ws.Range("D21").Value = '0.0₃0949'
$ws.Range("E21").Value = '  +2.08%  '

$ws.Range("E22").Value = '  -2.39%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '72.04'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +2.53%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '264.92'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +8.26%  '

$ws.Range("E25").Value = '  +1.00%  '

$ws.Range("E26").Value = '  -0.25%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '26.11'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.14%  '

$ws.Range("E28").Value = '  +1.14%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.78%  '

$ws.Range("E30").Value = '  +3.94%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '36.03'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +3.34%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '49.79'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.88%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '19.81'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.07%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '5.42'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.76%  '

$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("E36").Value = '  +1.44%  '

$ws.Range("E37").Value = '  +2.24%  '

$ws.Range("E38").Value = '  +1.96%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.00'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.41%  '

$ws.Range("E40").Value = '  +0.82%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '120.91'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.66%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '22.15'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.62%  '

$ws.Range("E43").Value = '  -0.86%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0302'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.87%  '

$ws.Range("D45").Value = '2.020.49'
$ws.Range("E45").Value = '  +2.25%  '

$ws.Range("E46").Value = '  +5.44%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.91'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +8.90%  '

$ws.Range("E48").Value = '  +2.14%  '

$ws.Range("E49").Value = '  +1.20%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '5.23'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +2.30%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '79.25'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +3.44%  '

Write-Output "Updated cryptos list"